# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and
# apply them to the runs that received a <w:rPr><w:rStyle .../></w:rPr>
# in the commit.

$d = $word.ActiveDocument

# --- 1. Create the character styles -----------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)          # wdStyleTypeCharacter
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)  # wdStyleTypeCharacter
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)          # wdStyleTypeCharacter
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608                     # wdColorNavy (BGR 0x800000 -> RGB 000080)
$gaNLinks.Font.Underline = 1                        # wdUnderlineSingle

# --- 2. Apply GaNStyle to all 4 "2022: Daty kampanii..." runs ---------

$dateText = "2022: Daty kampanii używające Gwiazdozbiór Bliźniąt: 14-23 lutego, 14-24 marca"

$rng = $d.Content
$rng.Start = 0
while ($rng.Find.Execute($dateText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# --- 3. Apply GaNParagraph to the long intro paragraph -----------------

$paragraphText = "Uczestniczysz w ogólnoświatowym przedsięwzięciu, którego celem jest obserwacja i odnotowanie najsłabszych widocznych gwiazd w celu zmierzenia zanieczyszczenia światłem w danym miejscu. Poprzez zlokalizowanie i obserwację  Gwiazdozbiór Bliźniąt na nocnym niebie oraz porównanie go do map nieba ludzie z całego świata będą mogli dowiedzieć się jaki wkład światło emitowane przez ich społeczność wnosi do  zanieczyszczenia światłem. To co dodasz do internetowej bazy danych pomoże udokumentować widoczne nocne niebo."

$rng = $d.Content
$rng.Start = 0
if ($rng.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- 4. Apply GaNLinks to the CzechGlobe credit run ---------------------

$linksText = " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$rng.Start = 0
if ($rng.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}

Write-Output "Done"
